# Add new daily-tracking rows 409..417 to the "Total" sheet, mirroring the
# existing day-over-day formula pattern (columns C,D,I,J,K,L,M are formulas
# relative to the row above; A,B,E,F,G,H are entered values), then repoint
# the rolling "latest day" summary formulas in N2:U2 at the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# raw input values for the new rows: Date, Confirmed(B), Recovered(E), F, Deaths(G), H
$newRows = @(
    @{ R = 409; A = 44376; B = 12889; E = 7; F = 0; G = 319; H = 12839 },
    @{ R = 410; A = 44377; B = 12897; E = 6; F = 0; G = 320; H = 12854 },
    @{ R = 411; A = 44378; B = 12900; E = 7; F = 0; G = 320; H = 12860 },
    @{ R = 412; A = 44379; B = 12915; E = 7; F = 0; G = 320; H = 12870 },
    @{ R = 413; A = 44380; B = 12915; E = 9; F = 0; G = 320; H = 12870 },
    @{ R = 414; A = 44381; B = 12915; E = 8; F = 0; G = 320; H = 12870 },
    @{ R = 415; A = 44382; B = 12931; E = 7; F = 0; G = 320; H = 12882 },
    @{ R = 416; A = 44383; B = 12934; E = 8; F = 0; G = 321; H = 12888 },
    @{ R = 417; A = 44384; B = 12944; E = 5; F = 0; G = 322; H = 12895 }
)

foreach ($row in $newRows) {
    $r = $row.R
    $prev = $r - 1

    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 1).NumberFormat = "d-mmm"
    $ws.Cells.Item($r, 1).HorizontalAlignment = -4108

    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 2).HorizontalAlignment = -4108

    $ws.Cells.Item($r, 3).Formula = "=B$r-H$r-G$r"
    $ws.Cells.Item($r, 3).HorizontalAlignment = -4108

    $ws.Cells.Item($r, 4).Formula = "=B$r-H$r"
    $ws.Cells.Item($r, 4).HorizontalAlignment = -4108

    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 5).HorizontalAlignment = -4108

    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 6).HorizontalAlignment = -4108

    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 7).HorizontalAlignment = -4108

    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 8).HorizontalAlignment = -4108

    $ws.Cells.Item($r, 9).Formula = "=B$r-B$prev"
    $ws.Cells.Item($r, 9).HorizontalAlignment = -4108

    $ws.Cells.Item($r, 10).Formula = "=H$r-H$prev"
    $ws.Cells.Item($r, 10).HorizontalAlignment = -4108

    $ws.Cells.Item($r, 11).Formula = "=G$r-G$prev"
    $ws.Cells.Item($r, 11).HorizontalAlignment = -4108

    $ws.Cells.Item($r, 12).Formula = "=E$r-E$prev"
    $ws.Cells.Item($r, 12).HorizontalAlignment = -4108

    $ws.Cells.Item($r, 13).Formula = "=E$r+F$r"
    $ws.Cells.Item($r, 13).HorizontalAlignment = -4108
}

# Roll the "latest day" summary block (top row, columns N:U) forward from
# row 408 to the new last row, 416 -- mirrors the prior edit that rolled it
# from 407 to 408, etc.
$ws.Range("N2").Formula = "=D416"
$ws.Range("O2").Formula = "=E416"
$ws.Range("P2").Formula = "=F416"
$ws.Range("Q2").Formula = "=K416"
$ws.Range("R2").Formula = "=J416"
$ws.Range("S2").Formula = "=I416"
$ws.Range("T2").Formula = "=A416"
$ws.Range("U2").Formula = "=L416"

# Reflect the new scroll position / selection the author left the sheet in.
$aw = $excel.ActiveWindow
$aw.ScrollRow = 408
$aw.ScrollColumn = 1
$ws.Range("I417").Select()

$wb.Save()
